$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (renamed + reordered)
$ws.Range("A1").Value = "VAB"
$ws.Range("B1").Value = "Sector VAB descripción"
$ws.Range("C1").Value = "Código"
$ws.Range("D1").Value = "Comarca"
$ws.Range("E1").Value = "Sector VAB código"
$ws.Range("F1").Value = "Año"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:vab"
$ws.Range("B2").Value = "iaest-measure:sector-vab-descripcion"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:comarca"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refPeriod"

# Row 3 - medida/dim classifiers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"

# Row 4 - datatypes
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "xsd:date"

# Row 5 - extra mapping file reference moves from C (ano) to F (Año)
$ws.Range("C5").Clear()
$ws.Range("A4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = "mapping-ano.xlsx"
